$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$recipeName = "매콤 두부 가지볶음"

$steps = @(
    "두부는 키친타올로 물기를 제거한 뒤 깍둑썰기 한다.",
    "가지는 반으로 자른 후 어슷하게 썬다.",
    "팬에 들기름을 두르고 마늘을 볶아 향을 낸다.",
    "두부와 가지를 넣고 중불에서 볶는다.",
    "간장, 고춧가루, 물을 넣고 뚜껑을 덮은 후 약불에서 3~4분간 졸인다.",
    "불을 끄고 쪽파를 넣어 마무리한다."
)

$ingredients = @(
    "두부",
    "가지",
    "다진 마늘",
    "간장",
    "들기름",
    "다진 쪽파",
    "고춧가루",
    "물"
)

# Column A: recipe name for rows 2-9
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $recipeName
}

# Column C: cooking steps for rows 2-7 (rows 8,9 have no step)
for ($i = 0; $i -lt $steps.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $steps[$i]
}

# Column B: ingredients for rows 2-9
for ($i = 0; $i -lt $ingredients.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $ingredients[$i]
}

$ws.Range("A1:B9").EntireColumn.AutoFit()
